# Auto-generated script applying scheduled-runner market data refresh
# to the Cactuar_Profits workbook. For each affected sheet/row/column this
# sets the refreshed numeric value, or clears the cell when the refreshed
# row no longer has a value in that column (diff removes the <c> element).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2185.6365  # H19: was 2301.6843
$ws.Cells.Item(19, 9).Value = 2093.1765  # I19: was 2227.3572
$ws.Cells.Item(19, 10).Value = 2500  # J19: was 2509.8
$ws.Cells.Item(19, 11).Value = 2093.1765  # K19: was 2227.3572
$ws.Cells.Item(19, 12).Value = 2500  # L19: was 2509.8
$ws.Cells.Item(19, 13).Value = -1918.1765  # M19: was -2052.3572
$ws.Cells.Item(19, 14).Value = -2850  # N19: was -2859.8
$ws.Cells.Item(51, 8).Value = 10617.407  # H51: was 10645.357
$ws.Cells.Item(51, 9).Value = 9232.666999999999  # I51: was 9699.333000000001
$ws.Cells.Item(51, 10).Value = 10790.5  # J51: was 10758.88
$ws.Cells.Item(51, 11).Value = 9232.666999999999  # K51: was 9699.333000000001
$ws.Cells.Item(51, 12).Value = 10790.5  # L51: was 10758.88
$ws.Cells.Item(51, 13).Value = -8748.666999999999  # M51: was -9215.333000000001
$ws.Cells.Item(51, 14).Value = -11758.5  # N51: was -11726.88
$ws.Cells.Item(76, 8).Value = 3246.625  # H76: was 3079.0908
$ws.Cells.Item(76, 9).Value = 2872.25  # I76: was 2564.5
$ws.Cells.Item(76, 10).Value = 3621  # J76: was 3696.6
$ws.Cells.Item(76, 11).Value = 2872.25  # K76: was 2564.5
$ws.Cells.Item(76, 12).Value = 3621  # L76: was 3696.6
$ws.Cells.Item(76, 13).Value = -2557.25  # M76: was -2249.5
$ws.Cells.Item(76, 14).Value = -4251  # N76: was -4326.6
$ws.Cells.Item(79, 8).Value = 3246.625  # H79: was 3079.0908
$ws.Cells.Item(79, 9).Value = 2872.25  # I79: was 2564.5
$ws.Cells.Item(79, 10).Value = 3621  # J79: was 3696.6
$ws.Cells.Item(79, 11).Value = 2872.25  # K79: was 2564.5
$ws.Cells.Item(79, 12).Value = 3621  # L79: was 3696.6
$ws.Cells.Item(79, 13).Value = -1780.25  # M79: was -1472.5
$ws.Cells.Item(79, 14).Value = -5805  # N79: was -5880.6
$ws.Cells.Item(101, 8).Value = 1612  # H101: was 1491.5
$ws.Cells.Item(101, 9).Value = 432.6  # I101: was 427
$ws.Cells.Item(101, 10).Value = 3086.25  # J101: was 3088.25
$ws.Cells.Item(101, 11).Value = 1297.8  # K101: was 1281
$ws.Cells.Item(101, 12).Value = 9258.75  # L101: was 9264.75
$ws.Cells.Item(101, 13).Value = 324.1999999999998  # M101: was 341
$ws.Cells.Item(101, 14).Value = -12502.75  # N101: was -12508.75
$ws.Cells.Item(111, 8).Value = 7080.7856  # H111: was 7877.1665
$ws.Cells.Item(111, 9).Value = 4414.2  # I111: was 4942.125
$ws.Cells.Item(111, 11).Value = 13242.6  # K111: was 14826.375
$ws.Cells.Item(111, 13).Value = -10175.6  # M111: was -11759.375
$ws.Cells.Item(137, 8).Value = 18072486  # H137: was 18072526
$ws.Cells.Item(137, 9).Value = 770800.4  # I137: was 834971.25
$ws.Cells.Item(137, 10).Value = 55559476  # J137: was 47622620
$ws.Cells.Item(137, 11).Value = 2312401.2  # K137: was 2504913.75
$ws.Cells.Item(137, 12).Value = 166678428  # L137: was 142867860
$ws.Cells.Item(137, 13).Value = -2309851.2  # M137: was -2502363.75
$ws.Cells.Item(137, 14).Value = -166683528  # N137: was -142872960

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 24816.73  # H32: was 26207.03
$ws.Cells.Item(32, 9).Value = 24552.182  # I32: was 26104.807
$ws.Cells.Item(32, 11).Value = 24552.182  # K32: was 26104.807
$ws.Cells.Item(32, 13).Value = -24265.182  # M32: was -25817.807
$ws.Cells.Item(43, 8).Value = 20000  # H43: was 0
$ws.Cells.Item(43, 10).Value = 20000  # J43: was 0
$ws.Cells.Item(43, 12).Value = 20000  # L43: was 0
$ws.Cells.Item(43, 14).Value = -20626  # N43: was None
$ws.Cells.Item(74, 8).Value = 1232.0435  # H74: was 1259.6666
$ws.Cells.Item(74, 9).Value = 1196.9  # I74: was 1225.2222
$ws.Cells.Item(74, 11).Value = 1196.9  # K74: was 1225.2222
$ws.Cells.Item(74, 13).Value = -322.9000000000001  # M74: was -351.2221999999999
$ws.Cells.Item(77, 8).Value = 1232.0435  # H77: was 1259.6666
$ws.Cells.Item(77, 9).Value = 1196.9  # I77: was 1225.2222
$ws.Cells.Item(77, 11).Value = 5984.5  # K77: was 6126.111
$ws.Cells.Item(77, 13).Value = -1616.5  # M77: was -1758.111
$ws.Cells.Item(102, 8).Value = 3997.6667  # H102: was 3999
$ws.Cells.Item(102, 10).Value = 3995  # J102: was 0
$ws.Cells.Item(102, 12).Value = 3995  # L102: was 0
$ws.Cells.Item(102, 14).Value = -7239  # N102: was None
$ws.Cells.Item(110, 8).Value = 1403.8462  # H110: was 1419.9231
$ws.Cells.Item(110, 9).Value = 1403.8462  # I110: was 1419.9231
$ws.Cells.Item(110, 11).Value = 1403.8462  # K110: was 1419.9231
$ws.Cells.Item(110, 13).Value = 641.1538  # M110: was 625.0769
$ws.Cells.Item(122, 8).Value = 3593.1943  # H122: was 3770.697
$ws.Cells.Item(122, 9).Value = 3593.1943  # I122: was 3835.4062
$ws.Cells.Item(122, 10).Value = 0  # J122: was 1700
$ws.Cells.Item(122, 11).Value = 10779.5829  # K122: was 11506.2186
$ws.Cells.Item(122, 12).Value = 0  # L122: was 5100
$ws.Cells.Item(122, 13).Value = -8329.582900000001  # M122: was -9056.2186
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -10000

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2133.889  # H107: was 2141.6296
$ws.Cells.Item(107, 9).Value = 1904.75  # I107: was 1917.8125
$ws.Cells.Item(107, 11).Value = 1904.75  # K107: was 1917.8125
$ws.Cells.Item(107, 13).Value = 15.25  # M107: was 2.1875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1122.6428  # H16: was 1194
$ws.Cells.Item(16, 9).Value = 775.1818  # I16: was 833.2
$ws.Cells.Item(16, 11).Value = 775.1818  # K16: was 833.2
$ws.Cells.Item(16, 13).Value = -488.1818  # M16: was -546.2
$ws.Cells.Item(22, 8).Value = 1429.2727  # H22: was 1098
$ws.Cells.Item(22, 9).Value = 1215.125  # I22: was 196
$ws.Cells.Item(22, 10).Value = 2000.3334  # J22: was 2000
$ws.Cells.Item(22, 11).Value = 1215.125  # K22: was 196
$ws.Cells.Item(22, 12).Value = 2000.3334  # L22: was 2000
$ws.Cells.Item(22, 13).Value = -865.125  # M22: was 154
$ws.Cells.Item(22, 14).Value = -2700.3334  # N22: was -2700
$ws.Cells.Item(31, 8).Value = 2713.4902  # H31: was 2680.5962
$ws.Cells.Item(31, 9).Value = 2284.8408  # I31: was 2284.9546
$ws.Cells.Item(31, 10).Value = 5407.857  # J31: was 4856.625
$ws.Cells.Item(31, 11).Value = 2284.8408  # K31: was 2284.9546
$ws.Cells.Item(31, 12).Value = 5407.857  # L31: was 4856.625
$ws.Cells.Item(31, 13).Value = -1989.8408  # M31: was -1989.9546
$ws.Cells.Item(31, 14).Value = -5997.857  # N31: was -5446.625
$ws.Cells.Item(34, 8).Value = 2713.4902  # H34: was 2680.5962
$ws.Cells.Item(34, 9).Value = 2284.8408  # I34: was 2284.9546
$ws.Cells.Item(34, 10).Value = 5407.857  # J34: was 4856.625
$ws.Cells.Item(34, 11).Value = 2284.8408  # K34: was 2284.9546
$ws.Cells.Item(34, 12).Value = 5407.857  # L34: was 4856.625
$ws.Cells.Item(34, 13).Value = -2082.8408  # M34: was -2082.9546
$ws.Cells.Item(34, 14).Value = -5811.857  # N34: was -5260.625
$ws.Cells.Item(58, 8).Value = 1795.7333  # H58: was 1567.2
$ws.Cells.Item(58, 9).Value = 1752.6428  # I58: was 1444
$ws.Cells.Item(58, 10).Value = 2399  # J58: was 2265.3333
$ws.Cells.Item(58, 11).Value = 1752.6428  # K58: was 1444
$ws.Cells.Item(58, 12).Value = 2399  # L58: was 2265.3333
$ws.Cells.Item(58, 13).Value = -1549.6428  # M58: was -1241
$ws.Cells.Item(58, 14).Value = -2805  # N58: was -2671.3333
$ws.Cells.Item(97, 8).Value = 23795  # H97: was 38975
$ws.Cells.Item(97, 10).Value = 23795  # J97: was 38975
$ws.Cells.Item(97, 12).Value = 23795  # L97: was 38975
$ws.Cells.Item(97, 14).Value = -25777  # N97: was -40957
$ws.Cells.Item(99, 8).Value = 10771.533  # H99: was 11405.286
$ws.Cells.Item(99, 9).Value = 14718.9  # I99: was 16143.223
$ws.Cells.Item(99, 10).Value = 2876.8  # J99: was 2877
$ws.Cells.Item(99, 11).Value = 14718.9  # K99: was 16143.223
$ws.Cells.Item(99, 12).Value = 2876.8  # L99: was 2877
$ws.Cells.Item(99, 13).Value = -13220.9  # M99: was -14645.223
$ws.Cells.Item(99, 14).Value = -5872.8  # N99: was -5873
$ws.Cells.Item(104, 8).Value = 60000  # H104: was 59999
$ws.Cells.Item(104, 9).Value = 0  # I104: was 59999
$ws.Cells.Item(104, 10).Value = 60000  # J104: was 0
$ws.Cells.Item(104, 11).Value = 0  # K104: was 59999
$ws.Cells.Item(104, 12).Value = 60000  # L104: was 0
$ws.Cells.Item(104, 13).ClearContents()  # M104: was -57378
$ws.Cells.Item(104, 14).Value = -65242  # N104: was None
$ws.Cells.Item(113, 8).Value = 1122.6428  # H113: was 1194
$ws.Cells.Item(113, 9).Value = 775.1818  # I113: was 833.2
$ws.Cells.Item(113, 11).Value = 775.1818  # K113: was 833.2
$ws.Cells.Item(113, 13).Value = 1394.8182  # M113: was 1336.8
$ws.Cells.Item(126, 8).Value = 10771.533  # H126: was 11405.286
$ws.Cells.Item(126, 9).Value = 14718.9  # I126: was 16143.223
$ws.Cells.Item(126, 10).Value = 2876.8  # J126: was 2877
$ws.Cells.Item(126, 11).Value = 44156.7  # K126: was 48429.669
$ws.Cells.Item(126, 12).Value = 8630.400000000001  # L126: was 8631
$ws.Cells.Item(126, 13).Value = -41686.7  # M126: was -45959.669
$ws.Cells.Item(126, 14).Value = -13570.4  # N126: was -13571
$ws.Cells.Item(136, 8).Value = 1795.7333  # H136: was 1567.2
$ws.Cells.Item(136, 9).Value = 1752.6428  # I136: was 1444
$ws.Cells.Item(136, 10).Value = 2399  # J136: was 2265.3333
$ws.Cells.Item(136, 11).Value = 5257.928400000001  # K136: was 4332
$ws.Cells.Item(136, 12).Value = 7197  # L136: was 6795.999899999999
$ws.Cells.Item(136, 13).Value = -2707.928400000001  # M136: was -1782
$ws.Cells.Item(136, 14).Value = -12297  # N136: was -11895.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 164001500  # H11: was 136667940
$ws.Cells.Item(11, 9).Value = 2516.3333  # I11: was 1912.25
$ws.Cells.Item(11, 11).Value = 7548.999899999999  # K11: was 5736.75
$ws.Cells.Item(11, 13).Value = -7408.999899999999  # M11: was -5596.75
$ws.Cells.Item(34, 8).Value = 1715.125  # H34: was 1421.85
$ws.Cells.Item(34, 9).Value = 525.1667  # I34: was 414.6
$ws.Cells.Item(34, 11).Value = 1575.5001  # K34: was 1243.8
$ws.Cells.Item(34, 13).Value = -1491.5001  # M34: was -1159.8
$ws.Cells.Item(120, 8).Value = 14999  # H120: was 14998.714
$ws.Cells.Item(120, 9).Value = 0  # I120: was 14997
$ws.Cells.Item(120, 11).Value = 0  # K120: was 44991
$ws.Cells.Item(120, 13).ClearContents()  # M120: was -40153
$ws.Cells.Item(128, 8).Value = 124899  # H128: was 124899.5
$ws.Cells.Item(128, 9).Value = 124899  # I128: was 124899.5
$ws.Cells.Item(128, 11).Value = 374697  # K128: was 374698.5
$ws.Cells.Item(128, 13).Value = -369717  # M128: was -369718.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 42030  # H49: was 0
$ws.Cells.Item(49, 10).Value = 42030  # J49: was 0
$ws.Cells.Item(49, 12).Value = 42030  # L49: was 0
$ws.Cells.Item(49, 14).Value = -42398  # N49: was None
$ws.Cells.Item(122, 8).Value = 316426.38  # H122: was 306864.5
$ws.Cells.Item(122, 9).Value = 668234.5600000001  # I122: was 626525.2
$ws.Cells.Item(122, 11).Value = 2004703.68  # K122: was 1879575.6
$ws.Cells.Item(122, 13).Value = -2002253.68  # M122: was -1877125.6
$ws.Cells.Item(123, 8).Value = 34432.4  # H123: was 34433.6
$ws.Cells.Item(123, 10).Value = 34432.4  # J123: was 34433.6
$ws.Cells.Item(123, 12).Value = 34432.4  # L123: was 34433.6
$ws.Cells.Item(123, 14).Value = -39332.4  # N123: was -39333.6
$ws.Cells.Item(132, 8).Value = 126930.5  # H132: was 144871.08
$ws.Cells.Item(132, 9).Value = 183181.55  # I132: was 201376.7
$ws.Cells.Item(132, 10).Value = 3178.2  # J132: was 3607
$ws.Cells.Item(132, 11).Value = 549544.6499999999  # K132: was 604130.1000000001
$ws.Cells.Item(132, 12).Value = 9534.599999999999  # L132: was 10821
$ws.Cells.Item(132, 13).Value = -547014.6499999999  # M132: was -601600.1000000001
$ws.Cells.Item(132, 14).Value = -14594.6  # N132: was -15881

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3951  # H40: was 4507.4165
$ws.Cells.Item(40, 9).Value = 2480  # I40: was 2763
$ws.Cells.Item(40, 11).Value = 2480  # K40: was 2763
$ws.Cells.Item(40, 13).Value = -2344  # M40: was -2627
$ws.Cells.Item(42, 8).Value = 0  # H42: was 30000
$ws.Cells.Item(42, 10).Value = 0  # J42: was 30000
$ws.Cells.Item(42, 12).Value = 0  # L42: was 30000
$ws.Cells.Item(42, 14).ClearContents()  # N42: was -31126
$ws.Cells.Item(49, 8).Value = 0  # H49: was 30000
$ws.Cells.Item(49, 10).Value = 0  # J49: was 30000
$ws.Cells.Item(49, 12).Value = 0  # L49: was 30000
$ws.Cells.Item(49, 14).ClearContents()  # N49: was -30294
$ws.Cells.Item(82, 8).Value = 2247.6316  # H82: was 2322.611
$ws.Cells.Item(82, 9).Value = 1648.4445  # I82: was 1742.25
$ws.Cells.Item(82, 11).Value = 1648.4445  # K82: was 1742.25
$ws.Cells.Item(82, 13).Value = -1287.4445  # M82: was -1381.25
$ws.Cells.Item(85, 8).Value = 2247.6316  # H85: was 2322.611
$ws.Cells.Item(85, 9).Value = 1648.4445  # I85: was 1742.25
$ws.Cells.Item(85, 11).Value = 1648.4445  # K85: was 1742.25
$ws.Cells.Item(85, 13).Value = -400.4445000000001  # M85: was -494.25
$ws.Cells.Item(93, 8).Value = 1501.5  # H93: was 1602.1666
$ws.Cells.Item(93, 9).Value = 1568.4166  # I93: was 1702.6
$ws.Cells.Item(93, 11).Value = 1568.4166  # K93: was 1702.6
$ws.Cells.Item(93, 13).Value = -320.4166  # M93: was -454.5999999999999
$ws.Cells.Item(122, 8).Value = 6532.7236  # H122: was 6819.952
$ws.Cells.Item(122, 9).Value = 3976.0334  # I122: was 4080.0386
$ws.Cells.Item(122, 10).Value = 11044.529  # J122: was 11272.3125
$ws.Cells.Item(122, 11).Value = 11928.1002  # K122: was 12240.1158
$ws.Cells.Item(122, 12).Value = 33133.587  # L122: was 33816.9375
$ws.Cells.Item(122, 13).Value = -9478.100199999999  # M122: was -9790.1158
$ws.Cells.Item(122, 14).Value = -38033.587  # N122: was -38716.9375
$ws.Cells.Item(135, 8).Value = 81250  # H135: was 80332.5
$ws.Cells.Item(135, 10).Value = 81250  # J135: was 80332.5
$ws.Cells.Item(135, 12).Value = 81250  # L135: was 80332.5
$ws.Cells.Item(135, 14).Value = -91390  # N135: was -90472.5
$ws.Cells.Item(136, 8).Value = 3848.2793  # H136: was 3872.8806
$ws.Cells.Item(136, 9).Value = 2772.7021  # I136: was 2785.152
$ws.Cells.Item(136, 11).Value = 8318.106299999999  # K136: was 8355.456
$ws.Cells.Item(136, 13).Value = -5768.106299999999  # M136: was -5805.456

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 8920.77  # H81: was 11019.444
$ws.Cells.Item(81, 9).Value = 6316.5  # I81: was 7079.8
$ws.Cells.Item(81, 10).Value = 11153  # J81: was 15944
$ws.Cells.Item(81, 11).Value = 12633  # K81: was 14159.6
$ws.Cells.Item(81, 12).Value = 22306  # L81: was 31888
$ws.Cells.Item(81, 13).Value = -11572  # M81: was -13098.6
$ws.Cells.Item(81, 14).Value = -24428  # N81: was -34010
$ws.Cells.Item(84, 8).Value = 8920.77  # H84: was 11019.444
$ws.Cells.Item(84, 9).Value = 6316.5  # I84: was 7079.8
$ws.Cells.Item(84, 10).Value = 11153  # J84: was 15944
$ws.Cells.Item(84, 11).Value = 63165  # K84: was 70798
$ws.Cells.Item(84, 12).Value = 111530  # L84: was 159440
$ws.Cells.Item(84, 13).Value = -57861  # M84: was -65494
$ws.Cells.Item(84, 14).Value = -122138  # N84: was -170048
$ws.Cells.Item(109, 8).Value = 0  # H109: was 49999
$ws.Cells.Item(109, 10).Value = 0  # J109: was 49999
$ws.Cells.Item(109, 12).Value = 0  # L109: was 49999
$ws.Cells.Item(109, 14).ClearContents()  # N109: was -52773
$ws.Cells.Item(126, 8).Value = 1888.4445  # H126: was 1963.25
$ws.Cells.Item(126, 9).Value = 1641.5714  # I126: was 1700.1666
$ws.Cells.Item(126, 11).Value = 4924.7142  # K126: was 5100.4998
$ws.Cells.Item(126, 13).Value = -2454.7142  # M126: was -2630.4998
$ws.Cells.Item(132, 8).Value = 1761.3418  # H132: was 1775.0769
$ws.Cells.Item(132, 9).Value = 1279.5555  # I132: was 1314.2354
$ws.Cells.Item(132, 11).Value = 3838.6665  # K132: was 3942.7062
$ws.Cells.Item(132, 13).Value = -1308.6665  # M132: was -1412.7062
$ws.Cells.Item(138, 8).Value = 98149.5  # H138: was 98150
$ws.Cells.Item(138, 10).Value = 98099  # J138: was 98100
$ws.Cells.Item(138, 12).Value = 98099  # L138: was 98100
$ws.Cells.Item(138, 14).Value = -108379  # N138: was -108380

